$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "D2" '27.869.83'
Set-TextValue "E2" '  -0.32%  '
Set-TextValue "D3" '1.630.04'
Set-TextValue "E3" '  -0.75%  '
Set-TextValue "E4" '  -0.20%  '
Set-TextValue "D5" '211.09'
Set-TextValue "E5" '  -0.75%  '
Set-TextValue "E6" '  -0.56%  '
Set-TextValue "D7" '0.998'
Set-TextValue "E7" '  -0.22%  '
Set-TextValue "E8" '  -0.66%  '
Set-TextValue "E9" '  -0.56%  '
Set-TextValue "E10" '  -0.44%  '
Set-TextValue "D11" '0.0879'
Set-TextValue "E11" '  -0.50%  '
Set-TextValue "D12" '1.861.74'
Set-TextValue "E12" '  -0.70%  '
Set-TextValue "D13" '1.623.77'
Set-TextValue "E13" '  -1.07%  '
Set-TextValue "E14" '  -1.52%  '
Set-TextValue "E15" '  -1.60%  '
Set-TextValue "E16" '  -0.41%  '
Set-TextValue "D17" '27.874.74'
Set-TextValue "E17" '  -0.26%  '
Set-TextValue "D18" '229.63'
Set-TextValue "E18" '  -1.60%  '
Set-TextValue "D19" '7.66'
Set-TextValue "E19" '  +0.54%  '
Set-TextValue "D20" '0.0₃0720'
Set-TextValue "E20" '  -0.34%  '
Set-TextValue "E21" '  -0.29%  '
Set-TextValue "D22" '4.34'
Set-TextValue "E22" '  -1.03%  '
Set-TextValue "E23" '  -5.03%  '
Set-TextValue "D24" '2.06'
Set-TextValue "E24" '  -1.67%  '
Set-TextValue "D25" '153.94'
Set-TextValue "E25" '  +0.43%  '
Set-TextValue "D26" '6.89'
Set-TextValue "E26" '  -0.12%  '
Set-TextValue "E27" '  -0.16%  '
Set-TextValue "D28" '15.53'
Set-TextValue "E28" '  -0.99%  '
Set-TextValue "E29" '  -0.17%  '
Set-TextValue "E30" '  -0.89%  '
Set-TextValue "E31" '  -0.74%  '
Set-TextValue "E32" '  -0.60%  '
Set-TextValue "D34" '1.393.76'
Set-TextValue "E34" '  -1.02%  '
Set-TextValue "E35" '  +0.47%  '
Set-TextValue "E36" '  +9.40%  '
Set-TextValue "E37" '  -1.13%  '
Set-TextValue "E39" '  -1.56%  '
Set-TextValue "D40" '0.852'
Set-TextValue "E40" '  -3.39%  '
Set-TextValue "B41" 'PaxDollar'
Set-TextValue "C41" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D41" '0.997'
Set-TextValue "E41" '  -0.25%  '
Set-TextValue "B42" 'WEMIXToken'
Set-TextValue "C42" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D42" '1.01'
Set-TextValue "E42" '  -1.53%  '
Set-TextValue "E43" '  -0.09%  '
Set-TextValue "D44" '65.82'
Set-TextValue "E44" '  -2.21%  '
Set-TextValue "D45" '5.43'
Set-TextValue "E45" '  -1.58%  '
Set-TextValue "D46" '1.768.57'
Set-TextValue "E46" '  -0.83%  '
Set-TextValue "E47" '  -2.82%  '
Set-TextValue "D48" '88.26'
Set-TextValue "E48" '  +0.26%  '
Set-TextValue "E49" '  +1.38%  '
Set-TextValue "E50" '  -0.48%  '
Set-TextValue "D51" '7.65'
Set-TextValue "E51" '  +0.39%  '
